$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.692.21'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '1.585.32'
$ws.Range("E3").Value = '  -2.48%  '
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").Value = '''207.13'
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("D6").Value = '''0.502'
$ws.Range("E6").Value = '  -3.17%  '
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("D8").Value = '''22.16'
$ws.Range("E8").Value = '  -4.46%  '
$ws.Range("D9").Value = '''0.252'
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = '''0.0591'
$ws.Range("E10").Value = '  -2.61%  '
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("D12").Value = '1.809.30'
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("D13").Value = '1.553.47'
$ws.Range("E13").Value = '  -4.47%  '
$ws.Range("E14").Value = '  -3.96%  '
$ws.Range("D15").Value = '''0.528'
$ws.Range("E15").Value = '  -4.75%  '
$ws.Range("D16").Value = '''63.52'
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '''219.63'
$ws.Range("E18").Value = '  -3.58%  '
$ws.Range("E19").Value = '  -3.12%  '
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").Value = '''4.13'
$ws.Range("E22").Value = '  -4.34%  '
$ws.Range("D23").Value = '''9.55'
$ws.Range("E23").Value = '  -3.56%  '
$ws.Range("D24").Value = '''1.96'
$ws.Range("E24").Value = '  -4.09%  '
$ws.Range("D25").Value = '''153.63'
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").Value = '''6.86'
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = '''15.10'
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("E29").Value = '  -4.30%  '
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").Value = '''0.0466'
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("D32").Value = '''3.22'
$ws.Range("E32").Value = '  -5.54%  '
$ws.Range("D33").Value = '1.362.73'
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").Value = '''2.94'
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("D35").Value = '''1.53'
$ws.Range("E35").Value = '  -4.06%  '
$ws.Range("D36").Value = '''0.972'
$ws.Range("E36").Value = '  -2.83%  '
$ws.Range("D37").Value = '''2.30'
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").Value = '''0.536'
$ws.Range("E39").Value = '  -3.02%  '
$ws.Range("D40").Value = '''0.821'
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("D42").Value = '''0.968'
$ws.Range("E42").Value = '  -2.86%  '
$ws.Range("D43").Value = '''63.74'
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").Value = '''2.17'
$ws.Range("E44").Value = '  +3.08%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '''1.73'
$ws.Range("E45").Value = '  -4.19%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''5.20'
$ws.Range("E46").Value = '  -3.68%  '
$ws.Range("D47").Value = '1.720.20'
$ws.Range("D48").Value = '''88.00'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("E49").Value = '  +10.79%  '
$ws.Range("D50").Value = '''0.0970'
$ws.Range("E50").Value = '  -3.80%  '
$ws.Range("E51").Value = '  -0.95%  '
